$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Amplicon")
$ws.Columns.Item(1).Insert()
$ws.Range("A1").Value = "amplicon_name"
$ws.Range("A2").Value = "test_amplicon"
$lo1 = $ws.ListObjects.Item(1)
$lo1.Unlist()
$lo2 = $ws.ListObjects.Item(1)
$lo2.Unlist()

$t = $ws.ListObjects.Add(1, $ws.Range("A1:N9"), $null, 1)
$t = $ws.ListObjects.Add(1, $ws.Range("O1:P9"), $null, 1)

$ws.ListObjects.Item("Table19").Name = "Table42"
$ws.ListObjects.Item("Table20").Name = "Table43"

$ws.ListObjects.Item("Table42").TableStyle = "TableStyleLight21"
$ws.ListObjects.Item("Table43").TableStyle = "TableStyleLight18"
